$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the entire row 62 (Lin Chun Lin's reservation); rows below shift up.
$ws.Rows.Item(62).Delete()
